$d = $word.ActiveDocument

$d.Content.Find.Execute("Ativação: 01/01/2014", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2021", 2)

$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11079086 - Herlandí de Souza Andrade", 2)

$d.Content.Find.Execute("Trabalhos em grupo; palestras e aulas expositivas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras", 2)

$d.Content.Find.Execute("Provas e trabalhos.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.", 2)

$d.Content.Find.Execute("Prova única com nota maior ou igual a 5,0 (cinco)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2)
